# Apply "normalization correction over dates": update the computed
# summary statistics (Sum, Average, Median, Max, Min) for rows 2-9
# to reflect a corrected normalization factor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 155.4430999755859
$ws.Range("C2").Value = 0.1166
$ws.Range("D2").Value = 0.1385000050067902
$ws.Range("E2").Value = 0.4257000088691711
$ws.Range("F2").Value = -0.1972000002861023

# Row 3
$ws.Range("B3").Value = 118.408203125
$ws.Range("C3").Value = 0.09470000000000001
$ws.Range("D3").Value = 0.0743
$ws.Range("E3").Value = 0.6294000148773193
$ws.Range("F3").Value = -0.1972000002861023

# Row 4
$ws.Range("B4").Value = 39.44770050048828
$ws.Range("C4").Value = 0.0461
$ws.Range("D4").Value = 0.0275
$ws.Range("E4").Value = 0.5091999769210815
$ws.Range("F4").Value = -0.1650999933481216

# Row 5
$ws.Range("B5").Value = -7.390900135040283
$ws.Range("C5").Value = -0.0078
$ws.Range("D5").Value = -0.0119
$ws.Range("E5").Value = 0.3248000144958496
$ws.Range("F5").Value = -0.2678999900817871

# Row 6
$ws.Range("B6").Value = -185.3421936035156
$ws.Range("C6").Value = -0.1466
$ws.Range("D6").Value = -0.1569
$ws.Range("E6").Value = 0.2806999981403351
$ws.Range("F6").Value = -0.3705999851226807

# Row 7
$ws.Range("B7").Value = -136.3917999267578
$ws.Range("C7").Value = -0.1374
$ws.Range("D7").Value = -0.1421999931335449
$ws.Range("E7").Value = 0.1283999979496002
$ws.Range("F7").Value = -0.3219999969005585

# Row 8
$ws.Range("B8").Value = 118.7403030395508
$ws.Range("C8").Value = 0.1064
$ws.Range("D8").Value = 0.1
$ws.Range("E8").Value = 0.3109999895095825
$ws.Range("F8").Value = -0.1312000006437302

# Row 9
$ws.Range("B9").Value = 102.914421081543
$ws.Range("C9").Value = 0.0139
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.6294000148773193
$ws.Range("F9").Value = -0.3705999851226807
